$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shift value Turno (A o B) from B to A
$ws.Range("C5").Value = "A"

# Update number of Full-Time workers
$ws.Range("C3").Value = 21

# Insert the "Viernes" row data and shift the Sabado/Domingo rows' values
# Row 10 (Miercoles): swap C10/E10
$ws.Range("C10").Value = 200
$ws.Range("E10").Value = 90

# Row 11 (Jueves)
$ws.Range("C11").Value = 200
$ws.Range("D11").Value = 50

# Row 12 - previously blank day label, now becomes "Viernes"
$ws.Range("B12").Value = "Viernes"
$ws.Range("C12").Value = 100
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = 100

# Row 13 - Sabado
$ws.Range("B13").Value = "Sábado"
$ws.Range("C13").Value = 400
$ws.Range("D13").Value = 200
$ws.Range("E13").Value = 500

# Row 14 - Domingo
$ws.Range("B14").Value = "Domingo"
$ws.Range("C14").Value = 400
$ws.Range("D14").Value = 400
$ws.Range("E14").Value = 100

# Update selection / active cell to G17
$ws.Range("G17").Select()

# Update window size
$excel.ActiveWindow.Width = 22188
$excel.ActiveWindow.Height = 9000
